$wb = $excel.ActiveWorkbook

# The "2022" sheet (physically sheet3.xml) gets a new entry appended.
$ws = $wb.Worksheets.Item("2022")

# Append the new row (row 5) with title / body(month) / image columns.
$ws.Range("A5").Value = "Selected for the Technology Innovation Hub - Cobotics Fellowship!"
$ws.Range("B5").Value = "July"
$ws.Range("C5").Value = "ihfc.PNG"

# Make "2022" the active sheet/tab, and set its selection to G8 - this also
# clears tabSelected on whichever sheet was previously active (2024).
$ws.Activate()
$ws.Range("G8").Select()
